# Qualitative-assessment workbook update:
#  - Fix scanner label for sub-9604 (CHUM -> Philips)
#  - Split sub-9611/9612/9613 (Oxford) into sessions of sub-9611: Ses1/Ses2/Ses3
#  - Split sub-9709/9710 (Milan) into sessions of sub-9709: Ses1/Ses2
#  - Split sub-10062/10063 (CEITEC) into sessions of sub-10062: Ses1/Ses2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename / relabel cells (Submission_ID, Image, Center columns) ---
$ws.Range("P6").Value = "Philips (9604)"
$ws.Range("P7").Value = "Philips (9604)"
$ws.Range("A10").Value = "sub-9611Ses1"
$ws.Range("B10").Value = "sub-9611Ses1_run-1_T2starw.nii.gz"
$ws.Range("P10").Value = "Oxford (9611Ses1)"
$ws.Range("A11").Value = "sub-9611Ses1"
$ws.Range("B11").Value = "sub-9611Ses1_run-2_T2starw.nii.gz"
$ws.Range("P11").Value = "Oxford (9611Ses1)"
$ws.Range("A12").Value = "sub-9611Ses2"
$ws.Range("B12").Value = "sub-9611Ses2_run-1_T2starw.nii.gz"
$ws.Range("P12").Value = "Oxford (9611Ses2)"
$ws.Range("A13").Value = "sub-9611Ses2"
$ws.Range("B13").Value = "sub-9611Ses2_run-2_T2starw.nii.gz"
$ws.Range("P13").Value = "Oxford (9611Ses2)"
$ws.Range("A14").Value = "sub-9611Ses3"
$ws.Range("B14").Value = "sub-9611Ses3_run-1_T2starw.nii.gz"
$ws.Range("P14").Value = "Oxford (9611Ses3)"
$ws.Range("A15").Value = "sub-9611Ses3"
$ws.Range("B15").Value = "sub-9611Ses3_run-2_T2starw.nii.gz"
$ws.Range("P15").Value = "Oxford (9611Ses3)"
$ws.Range("A18").Value = "sub-9709Ses1"
$ws.Range("B18").Value = "sub-9709Ses1_run-1_T2starw.nii.gz"
$ws.Range("P18").Value = "Milan (9709Ses1)"
$ws.Range("A19").Value = "sub-9709Ses1"
$ws.Range("B19").Value = "sub-9709Ses1_run-2_T2starw.nii.gz"
$ws.Range("P19").Value = "Milan (9709Ses1)"
$ws.Range("A20").Value = "sub-9709Ses2"
$ws.Range("B20").Value = "sub-9709Ses2_run-1_T2starw.nii.gz"
$ws.Range("P20").Value = "Milan (9709Ses2)"
$ws.Range("A21").Value = "sub-9709Ses2"
$ws.Range("B21").Value = "sub-9709Ses2_run-2_T2starw.nii.gz"
$ws.Range("P21").Value = "Milan (9709Ses2)"
$ws.Range("A22").Value = "sub-10062Ses1"
$ws.Range("B22").Value = "sub-10062Ses1_run-1_T2starw.nii.gz"
$ws.Range("P22").Value = "CEITEC (10062Ses1)"
$ws.Range("A23").Value = "sub-10062Ses1"
$ws.Range("B23").Value = "sub-10062Ses1_run-2_T2starw.nii.gz"
$ws.Range("P23").Value = "CEITEC (10062Ses1)"
$ws.Range("A24").Value = "sub-10062Ses2"
$ws.Range("B24").Value = "sub-10062Ses2_run-1_T2starw.nii.gz"
$ws.Range("P24").Value = "CEITEC (10062Ses2)"
$ws.Range("A25").Value = "sub-10062Ses2"
$ws.Range("B25").Value = "sub-10062Ses2_run-2_T2starw.nii.gz"
$ws.Range("P25").Value = "CEITEC (10062Ses2)"
$ws.Range("P32").Value = "Philips (9604)"
$ws.Range("P33").Value = "Philips (9604)"
$ws.Range("A36").Value = "sub-9611Ses1"
$ws.Range("B36").Value = "sub-9611Ses1_run-1_T2starw.nii.gz"
$ws.Range("P36").Value = "Oxford (9611Ses1)"
$ws.Range("A37").Value = "sub-9611Ses1"
$ws.Range("B37").Value = "sub-9611Ses1_run-2_T2starw.nii.gz"
$ws.Range("P37").Value = "Oxford (9611Ses1)"
$ws.Range("A38").Value = "sub-9611Ses2"
$ws.Range("B38").Value = "sub-9611Ses2_run-1_T2starw.nii.gz"
$ws.Range("P38").Value = "Oxford (9611Ses2)"
$ws.Range("A39").Value = "sub-9611Ses2"
$ws.Range("B39").Value = "sub-9611Ses2_run-2_T2starw.nii.gz"
$ws.Range("P39").Value = "Oxford (9611Ses2)"
$ws.Range("A40").Value = "sub-9611Ses3"
$ws.Range("B40").Value = "sub-9611Ses3_run-1_T2starw.nii.gz"
$ws.Range("P40").Value = "Oxford (9611Ses3)"
$ws.Range("A41").Value = "sub-9611Ses3"
$ws.Range("B41").Value = "sub-9611Ses3_run-2_T2starw.nii.gz"
$ws.Range("P41").Value = "Oxford (9611Ses3)"
$ws.Range("A44").Value = "sub-9709Ses1"
$ws.Range("B44").Value = "sub-9709Ses1_run-1_T2starw.nii.gz"
$ws.Range("P44").Value = "Milan (9709Ses1)"
$ws.Range("A45").Value = "sub-9709Ses1"
$ws.Range("B45").Value = "sub-9709Ses1_run-2_T2starw.nii.gz"
$ws.Range("P45").Value = "Milan (9709Ses1)"
$ws.Range("A46").Value = "sub-9709Ses2"
$ws.Range("B46").Value = "sub-9709Ses2_run-1_T2starw.nii.gz"
$ws.Range("P46").Value = "Milan (9709Ses2)"
$ws.Range("A47").Value = "sub-9709Ses2"
$ws.Range("B47").Value = "sub-9709Ses2_run-2_T2starw.nii.gz"
$ws.Range("P47").Value = "Milan (9709Ses2)"
$ws.Range("A48").Value = "sub-10062Ses1"
$ws.Range("B48").Value = "sub-10062Ses1_run-1_T2starw.nii.gz"
$ws.Range("P48").Value = "CEITEC (10062Ses1)"
$ws.Range("A49").Value = "sub-10062Ses1"
$ws.Range("B49").Value = "sub-10062Ses1_run-2_T2starw.nii.gz"
$ws.Range("P49").Value = "CEITEC (10062Ses1)"
$ws.Range("A50").Value = "sub-10062Ses2"
$ws.Range("B50").Value = "sub-10062Ses2_run-1_T2starw.nii.gz"
$ws.Range("P50").Value = "CEITEC (10062Ses2)"
$ws.Range("A51").Value = "sub-10062Ses2"
$ws.Range("B51").Value = "sub-10062Ses2_run-2_T2starw.nii.gz"
$ws.Range("P51").Value = "CEITEC (10062Ses2)"

# --- Column width adjustments (content reflow after subject renaming) ---
# (Input values are pre-compensated for Excel's pixel-quantized column-width
# storage so the saved width lands as close as possible to the target.)
$ws.Columns.Item(2).ColumnWidth = 29.25
$ws.Columns.Item(3).ColumnWidth = 26.09
$ws.Columns.Item(4).ColumnWidth = 25.25
$ws.Columns.Item(5).ColumnWidth = 33.09
$ws.Columns.Item(13).ColumnWidth = 42.42
$ws.Columns.Item(16).ColumnWidth = 17.92

# --- View state: zoom back to 100% and move selection to O33 ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("O33").Select() | Out-Null
